# This script applies a permutation of the per-row data blocks (columns
# D, M, N, O, P, Q, S) across rows 2-9 of the active sheet, as described
# by the target diff. The identifying columns (A, B, C, E, F, G, H, I,
# J, K, L, R, T) remain unchanged for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (source row's original data block
# is copied into destination row).
$mapping = @{
    2 = 5
    3 = 9
    4 = 6
    5 = 2
    6 = 8
    7 = 4
    8 = 7
    9 = 3
}

$cols = @("D", "M", "N", "O", "P", "Q", "S")

# Snapshot the original values for the affected columns/rows before making
# any changes, since several rows are both sources and destinations.
$original = @{}
foreach ($row in 2..9) {
    $original[$row] = @{}
    foreach ($col in $cols) {
        $original[$row][$col] = $ws.Range("$col$row").Value2
    }
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $original[$srcRow][$col]
    }
}
